# Update "想去人数" (F column) figures across the 展览 (Exhibitions),
# 演出 (Performances) and 全部类型 (All types) sheets to match the
# refreshed scrape output ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibit.Range("F2").Value  = 37
$wsExhibit.Range("F3").Value  = 103
$wsExhibit.Range("F4").Value  = 1494
$wsExhibit.Range("F5").Value  = 204
$wsExhibit.Range("F7").Value  = 39
$wsExhibit.Range("F8").Value  = 9842
$wsExhibit.Range("F9").Value  = 170
$wsExhibit.Range("F10").Value = 117
$wsExhibit.Range("F12").Value = 187
$wsExhibit.Range("F13").Value = 373
$wsExhibit.Range("F14").Value = 6807
$wsExhibit.Range("F15").Value = 1083
$wsExhibit.Range("F16").Value = 630
$wsExhibit.Range("F18").Value = 195

# 演出 (sheet2)
$wsShow.Range("F2").Value = 2
$wsShow.Range("F3").Value = 551

# 全部类型 (sheet4)
$wsAll.Range("F2").Value  = 37
$wsAll.Range("F3").Value  = 103
$wsAll.Range("F4").Value  = 1494
$wsAll.Range("F5").Value  = 204
$wsAll.Range("F6").Value  = 2
$wsAll.Range("F8").Value  = 39
$wsAll.Range("F9").Value  = 551
$wsAll.Range("F11").Value = 9842
$wsAll.Range("F12").Value = 170
$wsAll.Range("F13").Value = 117
$wsAll.Range("F15").Value = 187
$wsAll.Range("F16").Value = 373
$wsAll.Range("F17").Value = 6807
$wsAll.Range("F18").Value = 1083
$wsAll.Range("F19").Value = 630
$wsAll.Range("F21").Value = 195

Write-Host "Updated F-column counts on 展览, 演出, 全部类型 sheets"
